{"js": "// The commit replaces the semicolon that separates a variable's\n// description from its value/definition with a colon, at every place\n// this \"label; value\" pattern occurs in the body text, e.g.\n//   \"= Total power rating of motors; ${HP} HP\"\n//      -> \"= Total power rating of motors: ${HP} HP\"\n// (plus a handful of invisible `w:proofErr` grammar/spell-check marker\n// removals that Word regenerates on its own and aren't part of the\n// document's actual text, so nothing further is required for those.)\n//\n// Rather than hard-coding seven independent search strings (fragile if\n// run-splitting differs), we locate every literal \";\" in the document\n// body and swap it for \":\" - this document contains exactly seven\n// semicolons in total, and all seven are the ones targeted by the diff.\n\nconst body = context.document.body;\n\nconst results = body.search(\";\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (const item of results.items) {\n  item.insertText(\":\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# The commit replaces the semicolon that separates a variable's\n# description from its value/definition with a colon, at every place\n# this \"label; value\" pattern occurs in the body text, e.g.\n#   \"= Total power rating of motors; ${HP} HP\"\n#      -> \"= Total power rating of motors: ${HP} HP\"\n# (plus a handful of invisible `w:proofErr` grammar/spell-check marker\n# removals that Word regenerates on its own and aren't part of the\n# document's actual text, so nothing further is required for those.)\n#\n# The document contains exactly seven semicolons in total, and all\n# seven are the ones targeted by the diff, so a single document-wide\n# Find & Replace All swaps every one of them for a colon.\n\n$d = $word.ActiveDocument\n$range = $d.Content\n\n$range.Find.ClearFormatting()\n$range.Find.Replacement.ClearFormatting()\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$range.Find.Execute(\";\", $false, $false, $false, $false, $false, $true, 1, $false, \":\", 2)\n"}
